# Apply updated "dSF" (column F) values per the repulled data set.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -3
    3  = 3
    4  = 5
    5  = 6
    6  = -2
    7  = -1
    8  = -3
    9  = -3
    10 = -2
    11 = -1
    12 = 4
    14 = 1
    15 = -1
    17 = 1
    18 = -2
    19 = 4
    20 = -3
    21 = -1
    23 = 1
    24 = 6
    25 = 3
    26 = 9
    27 = 4
    28 = 3
    29 = 1
    30 = -2
    31 = -2
    32 = 7
    33 = -2
    34 = 3
    35 = -2
    36 = 4
    37 = -2
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
